$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "64.680.91"
$ws.Range("E2").Value = "  -3.71%  "
$ws.Range("D3").Value = "3.407.33"
$ws.Range("E3").Value = "  -4.51%  "
$ws.Range("E4").Value = "  +0.20%  "
Set-TextValue $ws.Range("D5") "580.01"
$ws.Range("E5").Value = "  -4.83%  "
Set-TextValue $ws.Range("D6") "132.99"
$ws.Range("E6").Value = "  -9.21%  "
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.406.59"
$ws.Range("E8").Value = "  -4.56%  "
Set-TextValue $ws.Range("D9") "0.481"
$ws.Range("E9").Value = "  -6.49%  "
$ws.Range("E10").Value = "  -10.09%  "
Set-TextValue $ws.Range("D11") "7.00"
$ws.Range("E11").Value = "  -11.36%  "
$ws.Range("E12").Value = "  -10.26%  "
$ws.Range("D13").Value = "3.985.15"
$ws.Range("E13").Value = "  -4.54%  "
Set-TextValue $ws.Range("D14") "0.0000176"
$ws.Range("E14").Value = "  -10.43%  "
$ws.Range("E15").Value = "  -1.93%  "
$ws.Range("D16").Value = "3.412.24"
$ws.Range("E16").Value = "  -4.22%  "
Set-TextValue $ws.Range("D17") "25.92"
$ws.Range("E17").Value = "  -11.11%  "
$ws.Range("D18").Value = "64.630.05"
$ws.Range("E18").Value = "  -3.41%  "
Set-TextValue $ws.Range("D19") "9.34"
$ws.Range("E19").Value = "  -16.11%  "
Set-TextValue $ws.Range("D20") "5.65"
$ws.Range("E20").Value = "  -9.39%  "
Set-TextValue $ws.Range("D21") "13.38"
$ws.Range("E21").Value = "  -9.19%  "
Set-TextValue $ws.Range("D22") "379.36"
$ws.Range("E22").Value = "  -11.31%  "
$ws.Range("E23").Value = "  +0.01%  "
Set-TextValue $ws.Range("D24") "0.537"
$ws.Range("E24").Value = "  -10.58%  "
Set-TextValue $ws.Range("D25") "71.55"
$ws.Range("E25").Value = "  -8.09%  "
$ws.Range("D26").Value = "3.546.04"
$ws.Range("E26").Value = "  -4.44%  "
$ws.Range("E27").Value = "  -12.55%  "
$ws.Range("E28").Value = "  +0.70%  "
Set-TextValue $ws.Range("D29") "7.17"
$ws.Range("E29").Value = "  -11.59%  "
Set-TextValue $ws.Range("D30") "2.17"
$ws.Range("E30").Value = "  -12.92%  "
Set-TextValue $ws.Range("D31") "7.90"
$ws.Range("E31").Value = "  -13.11%  "
$ws.Range("D32").Value = "3.427.78"
$ws.Range("E32").Value = "  -4.15%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  -9.74%  "
Set-TextValue $ws.Range("D35") "22.77"
$ws.Range("E35").Value = "  -7.14%  "
Set-TextValue $ws.Range("D36") "170.19"
$ws.Range("E36").Value = "  -4.29%  "
$ws.Range("E37").Value = "  -13.87%  "
Set-TextValue $ws.Range("D38") "6.58"
$ws.Range("E38").Value = "  -14.91%  "
$ws.Range("E39").Value = "  -13.41%  "
$ws.Range("E40").Value = "  -14.87%  "
Set-TextValue $ws.Range("D41") "0.0754"
$ws.Range("E41").Value = "  -9.23%  "
Set-TextValue $ws.Range("D42") "0.799"
$ws.Range("E42").Value = "  -8.10%  "
$ws.Range("E43").Value = "  +0.34%  "
Set-TextValue $ws.Range("D44") "41.91"
$ws.Range("E44").Value = "  -8.08%  "
Set-TextValue $ws.Range("D45") "4.23"
$ws.Range("E45").Value = "  -16.22%  "
Set-TextValue $ws.Range("D46") "1.59"
$ws.Range("E46").Value = "  -11.73%  "
Set-TextValue $ws.Range("D47") "1.10"
$ws.Range("E47").Value = "  -2.87%  "
Set-TextValue $ws.Range("D48") "22.26"
$ws.Range("E48").Value = "  -6.83%  "
Set-TextValue $ws.Range("D49") "6.44"
$ws.Range("E49").Value = "  -10.20%  "
$ws.Range("D50").Value = "2.192.67"
$ws.Range("E50").Value = "  -6.34%  "
Set-TextValue $ws.Range("D51") "19.91"
$ws.Range("E51").Value = "  -10.98%  "
